# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.032.02'
$ws.Range('E2').Value = '  -0.87%  '

$ws.Range('D3').Value = '2.450.70'
$ws.Range('E3').Value = '  -3.51%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.10'
$ws.Range('E5').Value = '  -0.12%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.39'
$ws.Range('E6').Value = '  -2.34%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.565'
$ws.Range('E8').Value = '  +0.00%  '

$ws.Range('D9').Value = '2.456.41'
$ws.Range('E9').Value = '  -3.19%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0979'
$ws.Range('E10').Value = '  -0.37%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.98'
$ws.Range('E12').Value = '  -3.90%  '

$ws.Range('E13').Value = '  -2.60%  '

$ws.Range('D14').Value = '2.886.85'
$ws.Range('E14').Value = '  -3.37%  '

$ws.Range('D15').Value = '57.963.10'
$ws.Range('E15').Value = '  -0.92%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.65'
$ws.Range('E16').Value = '  -2.70%  '

$ws.Range('D18').Value = '2.454.78'
$ws.Range('E18').Value = '  -3.24%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.41'
$ws.Range('E19').Value = '  -2.68%  '

$ws.Range('E20').Value = '  -1.49%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '315.02'
$ws.Range('E21').Value = '  -2.66%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.13'
$ws.Range('E22').Value = '  -0.59%  '

$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.24'
$ws.Range('E24').Value = '  +0.11%  '

$ws.Range('E25').Value = '  -0.77%  '

$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '2.568.49'
$ws.Range('E26').Value = '  -2.71%  '

$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.16%  '

$ws.Range('E28').Value = '  -1.55%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.26'
$ws.Range('E29').Value = '  -1.83%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '173.79'
$ws.Range('E30').Value = '  +3.16%  '

$ws.Range('E31').Value = '  -2.05%  '

$ws.Range('E32').Value = '  -1.78%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.16'
$ws.Range('E33').Value = '  -2.46%  '

$ws.Range('E34').Value = '  -5.28%  '

$ws.Range('E35').Value = '  -0.02%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  +0.02%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '17.86'
$ws.Range('E37').Value = '  -2.30%  '

$ws.Range('E38').Value = '  -5.78%  '

$ws.Range('E39').Value = '  -3.78%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.31'
$ws.Range('E40').Value = '  -0.46%  '

$ws.Range('E41').Value = '  +4.41%  '

$ws.Range('E42').Value = '  -2.66%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.40'
$ws.Range('E43').Value = '  -1.84%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '261.04'
$ws.Range('E44').Value = '  -6.39%  '

$ws.Range('E45').Value = '  -3.10%  '

$ws.Range('E46').Value = '  -4.70%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0922'
$ws.Range('E47').Value = '  +0.47%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '122.06'
$ws.Range('E48').Value = '  -6.27%  '

$ws.Range('E49').Value = '  -1.71%  '

$ws.Range('E50').Value = '  -1.78%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.07'
$ws.Range('E51').Value = '  -4.15%  '
